$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Drop the hyperlink that currently lives on A55 (the INE source URL) and
#    normalise that cell back to the plain "source" (italic) style instead of
#    the blue/underlined "HyperLink" style.
#    Deleting + re-inserting the row is the reliable way to get the cell back
#    onto the same style index ("source") used by its neighbours, since
#    Excel's native row-insert inherits the format of the row above it.
# ---------------------------------------------------------------------------
$ws.Range("A55").Hyperlinks.Delete()

$urlText = $ws.Range("A55").Value2
$ws.Rows(55).Delete()
$ws.Rows(55).Insert()
$ws.Range("A55").Value = $urlText

# ---------------------------------------------------------------------------
# 2. Space the "Source:" citation block out with blank rows in between each
#    line (matching the re-flowed layout in the target workbook), pushing the
#    "AFDB" / citation rows further down the sheet.
#    Inserting a row always inherits the format of the row immediately above
#    it, and every row in this block now shares the same "source" style, so
#    each newly inserted blank row naturally picks up that same style.
# ---------------------------------------------------------------------------
$ws.Rows(54).Insert()
$ws.Rows(56).Insert()
$ws.Rows(58).Insert()

# ---------------------------------------------------------------------------
# 3. Swap out the AFDB citation text for the new IFC citation (same cell,
#    new footnote text).
# ---------------------------------------------------------------------------
$ws.Range("A63").Value = '"Developing SMEs through Business Linkages", International Finance Corporation (IFC), 2008, p. 24, available at http://commdev.org/files/2328_file_Developing_SMEs_Through_Business_Linkages.pdf'
